$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Q(-1)" naive-error column is being inserted right after the row
# label column (i.e. before the current column B), for data rows 2-16
# only (row 1, the Q0..Q9 header row, is left untouched).
#
# Shift existing values in columns B:J one column to the right (into
# C:K) for every data row, working from the rightmost column back to
# the left so values are not clobbered before they are copied. Any
# previous value that was in column K falls off the end of the table,
# matching the diff (old K column values are discarded).
for ($r = 2; $r -le 16; $r++) {
    for ($c = 11; $c -ge 3; $c--) {
        $srcVal = $ws.Cells.Item($r, $c - 1).Value()
        $ws.Cells.Item($r, $c).Value = $srcVal
    }
}

# Now populate the freed-up column B with the new naive-component values.
$newColB = @{
    2  = -1.025188112727922
    3  = 0.08364543516793629
    4  = -0.1538585523806955
    5  = 0.7495351060200912
    6  = 0.03849281619118239
    7  = -0.2590580299438133
    8  = 0.01855976243503714
    9  = 0.1467044301255134
    10 = -0.1819613811903656
    11 = 0.4718454808444464
    12 = -0.08594117411414147
    13 = -0.07695400962807622
    14 = -0.5068991247689255
    15 = 0.6215838649243215
    16 = -0.2766911554241067
}

foreach ($r in $newColB.Keys) {
    $ws.Cells.Item($r, 2).Value = $newColB[$r]
}
